$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in dates for existing rows 15-17
$ws.Range("D15").Value = " 13-02-2020"
$ws.Range("D16").Value = " 13-02-2021"
$ws.Range("D17").Value = " 13-02-2022"

# Billet 1 (row 13): probleme d'icon
$ws.Range("D13").Value = " 14-02-2018"
$ws.Range("E13").Value = "Icon"
$ws.Range("F13").Value = "L'icon de l'app qui apparait est l'icon par defaut"

# Billet 2 (row 14): probleme de chronometre
$ws.Range("F14").Value = "Chronomètre ne recommance pas a chaque fois qu'on entre un mot"
$ws.Range("E14").Value = "Temps"
$ws.Range("D14").Value = " 14-02-2019"

# Widen column F to better fit the new, longer descriptions
$ws.Columns.Item(6).ColumnWidth = 62.92

# Move the active selection, as left by the editor after the last edit
[void]$ws.Range("F19").Select()
